# Update kansensya_pcr.xlsx: refresh daily COVID-19 test/case figures and
# the "municipal outsiders" footnote, then leave the "kobe" sheet focused
# as the active tab/selection (matching the upstream data-refresh bot).

$wb = $excel.ActiveWorkbook

$wsAll   = $wb.Worksheets.Item("all")
$wsKobe  = $wb.Worksheets.Item("kobe")
$wsOther = $wb.Worksheets.Item("other")

# ---------------------------------------------------------------------
# "all" sheet: update the last summary row (row 34) and the footnote
# shared-string text in row 35.
# ---------------------------------------------------------------------
$wsAll.Range("C34").Value = 275
$wsAll.Range("D34").Value = 81
$wsAll.Range("E34").Value = 70

$footnote = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272・276・277例目（計16件）は市外在住者です。"
$wsAll.Range("B35").Value = $footnote

# ---------------------------------------------------------------------
# "kobe" sheet: refresh the trailing daily rows (79-89) and the matching
# footnote string in B90 (same text as above so the two collapse back to
# a single shared string).
# ---------------------------------------------------------------------
$wsKobe.Range("B79").Value = 142
$wsKobe.Range("C79").Value = 2092

$wsKobe.Range("C80").Value = 2125
$wsKobe.Range("C81").Value = 2183
$wsKobe.Range("C82").Value = 2232
$wsKobe.Range("C83").Value = 2269
$wsKobe.Range("C84").Value = 2341
$wsKobe.Range("C85").Value = 2407
$wsKobe.Range("C86").Value = 2470
$wsKobe.Range("C87").Value = 2552
$wsKobe.Range("C88").Value = 2590

$wsKobe.Range("B89").Value = 69
$wsKobe.Range("C89").Value = 2659
$wsKobe.Range("F89").Value = 76
$wsKobe.Range("G89").Value = 66

$wsKobe.Range("B90").Value = $footnote

# ---------------------------------------------------------------------
# View/selection bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------

# "all": move the remembered selection from A34 to A33.
$wsAll.Range("A33").Select()

# "other": move the remembered pane/selection one column to the right.
$wsOther.Range("J64").Select()

# "kobe": zoom to 85%, move the remembered selection to B89, and make it
# the active/selected sheet/tab (also updates workbookView.activeTab).
$wsKobe.Activate()
$excel.ActiveWindow.Zoom = 85
$wsKobe.Range("B89").Select()
